$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-12 from 45212 to 45221
$ws.Range("C2:C12").Value = 45221
